# Motor works good at 100 ms
# Add two new log entries (row 6 and a new row 7) to the time-tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: finish filling in the previously-empty row
#   Date: 2017-03-19, Time: 15:00, Duration: 5h, Task: "Поставил двигатель"
$ws.Range("A6").NumberFormat = $ws.Range("A5").NumberFormat
$ws.Range("A6").Value = 42813

$ws.Range("B6").NumberFormat = $ws.Range("B5").NumberFormat
$ws.Range("B6").Value = 0.625

$ws.Range("C6").Value = 5

$ws.Range("D6").Value = "Поставил двигатель"

# Row 7: new entry
#   Date: 2017-03-21, Time: 14:00, Task: "Тестирование двигателя"
$ws.Range("A7").NumberFormat = $ws.Range("A5").NumberFormat
$ws.Range("A7").Value = 42815

$ws.Range("B7").NumberFormat = $ws.Range("B5").NumberFormat
$ws.Range("B7").Value = 0.58333333333333337

$ws.Range("D7").Value = "Тестирование двигателя"

# Update the active selection to reflect where the user ended up
$ws.Range("D8").Select()
